# LOT2030.docx edit script
# Applies the content changes described by the commit diff:
#  - bump the "Ativacao" date
#  - re-shuffle / rewrite the objectives, professors, summary, program,
#    evaluation and bibliography paragraphs
#  - delete the two now-empty paragraphs

$d = $word.ActiveDocument

# Replace the first occurrence of $findText with $replaceText, but only
# search inside paragraph number $paraIndex (1-based, current/live
# numbering at the time this call runs) -- this keeps every edit local
# so that later paragraphs are never accidentally touched by an earlier
# step (and vice versa), even though several snippets of text get
# duplicated/moved around by this edit.
function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $rng = $d.Paragraphs.Item($paraIndex).Range
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find failed in paragraph $paraIndex for: $findText"
    }
}

# ---------------------------------------------------------------------
# 1) Créditos bullet list (paragraph 4): Ativação date
# ---------------------------------------------------------------------
Replace-InParagraph 4 "Ativação: 01/01/2025" "Ativação: 15/07/2025"

# ---------------------------------------------------------------------
# 2) Objetivos paragraph text (paragraph 6)
# ---------------------------------------------------------------------
Replace-InParagraph 6 `
    "Introdução à tecnologia de conversão de biomassa vegetal para estudantes de Engenharia Bioquímica, abordando os principais processos tecnológicos do setor e seus métodos de controle que incluem: celulose e papel; derivados de celulose; carvão vegetal e frações monoméricas por hidrólise." `
    "A disciplina aborda uma série de processos industriais que são utilizados no fracionamento e na conversão da biomassa vegetal ligninficada em produtos elaborados como celulose e papel, derivados de celulose, carvão e açúcares. Visita supervisionada prevista."

# ---------------------------------------------------------------------
# 3) Docente(s) Responsável(eis) bullet list (paragraph 9) -> becomes
#    the long "programa" description (two runs, first keeps its
#    trailing <w:br/>)
# ---------------------------------------------------------------------
Replace-InParagraph 9 "2143261 - André Luis Ferraz" `
    "Introdução à tecnologia de conversão de biomassa vegetal para estudantes de Engenharia Bioquímica, abordando os principais processos tecnológicos do setor e seus métodos de controle que incluem: celulose e papel; derivados de celulose; carvão vegetal e frações monoméricas por hidrólise."

Replace-InParagraph 9 "5111420 - Talita Martins Lacerda" `
    "1. Breve introdução sobre a disponibilidade da biomassa: tipos de biomassa lignificada, produção por reflorestamento, resíduos agrícolas, características celulares, composição química. 2. Produção de celuloses e papel: mercado mundial de celulose e papel, processos de polpação mecânica, kraft e sulfito; braqueamento de pastas celulósicas; recuperação de inorgânicos e geração de energia; métodos de controle de processo; características físico-químicas e métodos de produção de papel. 3. Produção de derivados de celulose: formação do celulosato em meio alcalino, nitrato de celulose, xantato de celulose e a produção de fibras têxteis de ""viscose"", acetato de celulose, carboximetil celulose, etil e propilcelulose, alongamento da cadeia celulósica com epóxidos. 4. Conversão térmica e produção de carvão vegetal: secagem da madeira e estabilização dimensional, processos termomecânicos e produção de aglomerados; energia de biomassa vegetal, queima para geração direta de energia; produção de carvão vegetal. 5. Produção de açúcares e derivados por hidrólise: hidrólise ácida e processos de pré-tratamento para desestruturação da parede celular. 6. Processos integrados para a conversão de biomassa: indústrias modernas que aplicam os conceitos de uso integrado da biomassa vegetal. 7. Visita supervisionada a laboratórios e indústrias, a depender da viabilidade no momento do oferecimento da disciplina."

# ---------------------------------------------------------------------
# 4) "Programa resumido" PT paragraph (paragraph 11) -> evaluation
#    method sentence
# ---------------------------------------------------------------------
Replace-InParagraph 11 `
    "Estrutura e ultraestrutura dos materiais lignocelulósicos, celulose, hemiceluloses e outras polioses. Lignina, extrativos e composição da casca. Reações em meio ácido, meio alcalino. A disciplina aborda uma série de processos industriais que são utilizados no fracionamento e na conversão da biomassa vegetal ligninficada em produtos elaborados como celulose e papel, derivados de celulose, carvão e açúcares." `
    "A avaliação será feita por meio de provas escritas."

# ---------------------------------------------------------------------
# 5) "Programa" long PT paragraph (paragraph 14) -> evaluation
#    criteria sentence
# ---------------------------------------------------------------------
Replace-InParagraph 14 `
    "1. Breve introdução sobre a disponibilidade da biomassa: tipos de biomassa lignificada, produção por reflorestamento, resíduos agrícolas, características celulares, composição química. 2. Produção de celuloses e papel: mercado mundial de celulose e papel, processos de polpação mecânica, kraft e sulfito; braqueamento de pastas celulósicas; recuperação de inorgânicos e geração de energia; métodos de controle de processo; características físico-químicas e métodos de produção de papel. 3. Produção de derivados de celulose: formação do celulosato em meio alcalino, nitrato de celulose, xantato de celulose e a produção de fibras têxteis de ""viscose"", acetato de celulose, carboximetil celulose, etil e propilcelulose, alongamento da cadeia celulósica com epóxidos. 4. Conversão térmica e produção de carvão vegetal: secagem da madeira e estabilização dimensional, processos termomecânicos e produção de aglomerados; energia de biomassa vegetal, queima para geração direta de energia; produção de carvão vegetal. 5. Produção de açúcares e derivados por hidrólise: hidrólise ácida e processos de pré-tratamento para desestruturação da parede celular. 6. Processos integrados para a conversão de biomassa: indústrias modernas que aplicam os conceitos de uso integrado da biomassa vegetal." `
    "A nota final (NF) será calculada da seguintes maneira: NF=(P1+P2)/2 x 0,9 + Estudo de Caso x 0,1."

# ---------------------------------------------------------------------
# 6) Italic EN "Programa" list (paragraph 15, was 6 items split
#    across <w:br/>) -> single run with 7 items, numbering normalised
#    ("1." -> "1. ")
# ---------------------------------------------------------------------
$p15 = $d.Paragraphs.Item(15)
$p15.Range.Text = "1. Brief introduction about the availability of biomass: kinds of lignified biomass, forestry, agricultural residues, cellular characteristics, chemical composition. 2. Production of pulp and paper: world market of pulp and paper, mechanical, kraft and sulfite pulping processes; pulp bleaching; recovery of inorganics and generation of energy; methods of process control, physical-chemical characteristics and methods of paper production. 3. Production of cellulose derivatives: alkali cellulose, cellulose nitrate, xanthate of cellulose and viscose production, cellulose acetate, carboxymethyl cellulose, ethyl and propyl cellulose, stretching of cellulosic chain with epoxides. 4. Thermo-conversion and vegetal charcoal production; wood drying and dimensional stabilization, thermomechanical processes and fiberboard production; energy from vegetal biomass. 5. Production of sugars and derivatives by acid hydrolysis and processes of pretreatment for cell wall disruption. 6. Integrated processes for biomass conversion: modern industries that apply the integrated use of vegetal biomass in the biorefinery concept. 7. Supervised visits to laboratories and industries, depending on feasibility at the time the discipline is offered"

# ---------------------------------------------------------------------
# 7) Avaliação bullet list (paragraph 17): swap the three answer runs
#    around, and bring in the (ex) bibliography text under "Critério:"
# ---------------------------------------------------------------------
Replace-InParagraph 17 "A avaliação será feita por meio de provas escritas." `
    "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada como MR=(NF=PR)/2."

Replace-InParagraph 17 "A nota final (NF) será calculada da seguintes maneira: NF=(P1+P2)/2 x 0,9 + Estudo de Caso x 0,1." `
    "1. EK, M., GELLERSTEDT, G., HENRIKSSON, G. Wood Chemistry and Wood Biotechnology (Vol 1) e Pulping Chemistry and Technology (Vol 2). Berlin: Walter de Gruyter, 2009.^l2. KLEMM, D., PHILIPP, B., HEINZE, T., HEINZE, U., WAGENKNECHT, U. Comprehensive Cellulose Chemistry (Volume 2-Functionalization of Cellulose). Berlin: Wyley, 1998.^l3. FENGEL, D., WEGENER, G. Wood Chemistry, Ultrastruture, Reactions. Berlin: Walter de Gruyter,1989.^l"

Replace-InParagraph 17 "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada como MR=(NF=PR)/2." `
    "2143261 - André Luis Ferraz"

# ---------------------------------------------------------------------
# 8) Bibliography paragraph (paragraph 19) -> now just the 2nd
#    professor's line
# ---------------------------------------------------------------------
Replace-InParagraph 19 `
    "1. EK, M., GELLERSTEDT, G., HENRIKSSON, G. Wood Chemistry and Wood Biotechnology (Vol 1) e Pulping Chemistry and Technology (Vol 2). Berlin: Walter de Gruyter, 2009." `
    "5111420 - Talita Martins Lacerda"

# ---------------------------------------------------------------------
# 9) Delete the two paragraphs left empty / orphaned by the moves
#    above (process from the higher index down so earlier indices
#    stay valid for the second delete)
# ---------------------------------------------------------------------
$d.Paragraphs.Item(12).Range.Delete()  # italic EN "The discipline deals..."
$d.Paragraphs.Item(7).Range.Delete()   # empty italic run after Objetivos
